# Apply updated dSF (column F) values as per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    4  = -3
    5  = -1
    8  = 3
    10 = -2
    12 = 3
    18 = -3
    19 = 2
    20 = 4
    21 = -4
    22 = -2
    23 = -4
    24 = -5
    25 = -1
    26 = -1
    27 = 2
    28 = -1
    29 = 5
    31 = 11
    32 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
